# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" (column H) for the
# 787154a6-d431-4486-ba60-c634e9d6534b.md row (row 6) on the per-locale
# report sheets, reflecting a fresh handoff xliff generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("H6").Value = "2016-08-17 16:41:27"
$wsDeDe.Range("H6").Value = "2016-08-17 16:41:32"
